# Generate Report for Archive
# - Status moves from "Ready for handoff" to "In Translation" for the two
#   rows in every sheet (Overview uses E/F per-locale status columns,
#   zh-cn / de-de each use their own Status column C).
# - The Status column(s) are narrowed to fit the new, shorter text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: columns E (zh-cn) and F (de-de) hold the status text ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F3").Value = "In Translation"
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet: column C is "Status" ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2:C3").Value = "In Translation"
$zhcn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet: column C is "Status" ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2:C3").Value = "In Translation"
$dede.Columns.Item(3).ColumnWidth = 12.5
